$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from an existing header
# cell onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record columns: every player row gets the team's W-L-T record.
$ws.Range("AD2:AD49").Value = 84
$ws.Range("AE2:AE49").Value = 78
$ws.Range("AF2:AF49").Value = 0
